# Replace the scraped product rows (2-11) with the new 2023-06-20 scrape
# results. Only cell text is touched (existing hyperlinks on column C stay
# bound to their original target URLs - the workbook's rels are untouched
# by this edit, only the displayed/shared-string text changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = '(IBOX) Apple iPhone 14 512GB 256GB 128GB 6.1" inch Resmi Indo TAM'
$ws.Cells.Item(2, 2).Value = 'Rp13.580.000'
$ws.Cells.Item(2, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2022/10/21/f4c2f823-4297-4da7-a416-dd7658ccaf91.jpg'
$ws.Cells.Item(2, 4).Value = 'Putra Group'
$ws.Cells.Item(2, 5).Value = 'Jakarta Pusat'
$ws.Cells.Item(2, 6).Value = 'Terjual 100+'
$ws.Cells.Item(2, 7).Value = 'Tokopedia'

# Row 3
$ws.Cells.Item(3, 1).Value = 'iPhone 14 Promax Garansi Resmi'
$ws.Cells.Item(3, 2).Value = 'Rp29.000.000'
$ws.Cells.Item(3, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2022/12/17/deb74364-7a7a-44a7-a5db-b0ca4332dd8e.png'
$ws.Cells.Item(3, 4).Value = 'PT Pratama Sntra Semesta'
$ws.Cells.Item(3, 5).Value = 'Jakarta Barat'
$ws.Cells.Item(3, 6).Value = 'Terjual 500+'
$ws.Cells.Item(3, 7).Value = 'Tokopedia'

# Row 4
$ws.Cells.Item(4, 1).Value = '[PASTI RESMI] Apple iPhone 14 PRO MAX 1TB 512GB 256GB 128GB Resmi'
$ws.Cells.Item(4, 2).Value = 'Rp20.105.000'
$ws.Cells.Item(4, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/2/17/237017e6-ae9d-44bc-b4a9-a5a315ec563c.png'
$ws.Cells.Item(4, 4).Value = 'Apple Bank'
$ws.Cells.Item(4, 5).Value = 'Jakarta Barat'
$ws.Cells.Item(4, 6).Value = 'Terjual 30+'
$ws.Cells.Item(4, 7).Value = 'Tokopedia'

# Row 5
$ws.Cells.Item(5, 1).Value = 'RESMI iPhone 14 5G 128 / 256 / 512 Midnight Starlight Purple Blue - 128GB SINGLE, RED'
$ws.Cells.Item(5, 2).Value = 'Rp12.319.000'
$ws.Cells.Item(5, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2022/9/15/34d4b699-5a4a-46d6-9220-4593e02297a2.jpg'
$ws.Cells.Item(5, 4).Value = 'tokohapedia'
$ws.Cells.Item(5, 5).Value = 'Jakarta Pusat'
$ws.Cells.Item(5, 6).Value = 'Terjual 50+'
$ws.Cells.Item(5, 7).Value = 'Tokopedia'

# Row 6 (note: the "sold" column has no value for this listing)
$ws.Cells.Item(6, 1).Value = 'Iphone 14 128 New garansi resmi blm aktif'
$ws.Cells.Item(6, 2).Value = 'Rp13.700.000'
$ws.Cells.Item(6, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/6/18/4d428a92-f853-4fb4-a207-963f7e90e3aa.jpg'
$ws.Cells.Item(6, 4).Value = 'Ilham Pusat HP Second BPP'
$ws.Cells.Item(6, 5).Value = 'Balikpapan'
$ws.Cells.Item(6, 6).ClearContents()
$ws.Cells.Item(6, 7).Value = 'Tokopedia'

# Row 7
$ws.Cells.Item(7, 1).Value = 'Apple iPhone 14 Pro Garansi Resmi - 128GB 256GB 512GB 1TB'
$ws.Cells.Item(7, 2).Value = 'Rp17.269.000'
$ws.Cells.Item(7, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/2/4/6a115d00-bc7d-4d29-beb1-7939b44983f5.jpg'
$ws.Cells.Item(7, 4).Value = 'iSmile Official Store'
$ws.Cells.Item(7, 5).Value = 'Jakarta Pusat'
$ws.Cells.Item(7, 6).Value = 'Terjual 750+'
$ws.Cells.Item(7, 7).Value = 'Tokopedia'

# Row 8
$ws.Cells.Item(8, 1).Value = 'iPhone 14 Garansi Resmi'
$ws.Cells.Item(8, 2).Value = 'Rp13.390.000'
$ws.Cells.Item(8, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/3/9/2cc67a4b-0bd9-4ae8-a35b-1d1bdec850ec.png'
$ws.Cells.Item(8, 4).Value = 'PT Pratama Sntra Semesta'
$ws.Cells.Item(8, 5).Value = 'Jakarta Barat'
$ws.Cells.Item(8, 6).Value = 'Terjual 250+'
$ws.Cells.Item(8, 7).Value = 'Tokopedia'

# Row 9
$ws.Cells.Item(9, 1).Value = 'iPhone 14 Pro Garansi Resmi'
$ws.Cells.Item(9, 2).Value = 'Rp17.910.000'
$ws.Cells.Item(9, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2022/12/17/fd8a942f-df0a-4c9a-9cef-014f16d32bc6.png'
$ws.Cells.Item(9, 4).Value = 'PT Pratama Sntra Semesta'
$ws.Cells.Item(9, 5).Value = 'Jakarta Barat'
$ws.Cells.Item(9, 6).Value = 'Terjual 500+'
$ws.Cells.Item(9, 7).Value = 'Tokopedia'

# Row 10
$ws.Cells.Item(10, 1).Value = 'Apple iPhone 14 Pro Max Garansi Resmi - 128GB 256GB 512GB 1TB Promax'
$ws.Cells.Item(10, 2).Value = 'Rp19.159.000'
$ws.Cells.Item(10, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/2/15/5e89831a-e8a9-4f20-af9c-c701c347d1c3.jpg'
$ws.Cells.Item(10, 4).Value = 'iSmile Official Store'
$ws.Cells.Item(10, 5).Value = 'Jakarta Pusat'
$ws.Cells.Item(10, 6).Value = 'Terjual 500+'
$ws.Cells.Item(10, 7).Value = 'Tokopedia'

# Row 11
$ws.Cells.Item(11, 1).Value = 'Apple iPhone 14 Garansi Resmi - 128GB 256GB 512GB'
$ws.Cells.Item(11, 2).Value = 'Rp13.389.000'
$ws.Cells.Item(11, 3).Value = 'https://images.tokopedia.net/img/cache/200-square/VqbcmM/2023/2/15/0c4eabda-32af-4afa-b5b4-5a1363c705c4.jpg'
$ws.Cells.Item(11, 4).Value = 'iSmile Official Store'
$ws.Cells.Item(11, 5).Value = 'Jakarta Pusat'
$ws.Cells.Item(11, 6).Value = 'Terjual 250+'
$ws.Cells.Item(11, 7).Value = 'Tokopedia'
